$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 378, pushing the existing data (rows 378..455)
# down to (379..456), growing the used range to A1:R456.
$ws.Rows.Item(378).Insert()

# Populate the newly inserted row 378 with the new weekly record.
$ws.Cells.Item(378, 1).Value2 = 3
$ws.Cells.Item(378, 2).Value2 = "Femacal de La Calera"
$ws.Cells.Item(378, 3).Value2 = "Coquimbo"
$ws.Cells.Item(378, 4).Value2 = 44694
$ws.Cells.Item(378, 5).Value2 = 5
$ws.Cells.Item(378, 6).Value2 = 100112021
$ws.Cells.Item(378, 7).Value2 = "Ají"
$ws.Cells.Item(378, 8).Value2 = "Americana (o)"
$ws.Cells.Item(378, 9).Value2 = "Primera"
$ws.Cells.Item(378, 10).Value2 = 73
$ws.Cells.Item(378, 11).Value2 = 26000
$ws.Cells.Item(378, 12).Value2 = 27000
$ws.Cells.Item(378, 13).Value2 = 26479
$ws.Cells.Item(378, 14).Value2 = "$/caja 15 kilos"
$ws.Cells.Item(378, 15).Value2 = "Región de Arica y Parinacota"
$ws.Cells.Item(378, 16).Value2 = 1765
$ws.Cells.Item(378, 17).Value2 = 15
$ws.Cells.Item(378, 18).Value2 = "Hortaliza"
